# Update "想去人数" (number interested) counts on the 展览 (Exhibitions)
# and 全部类型 (All Types) sheets. Each listed cell's value is incremented
# as published by the latest generated-data refresh (gh-pages output).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 320
$ws1.Range("F5").Value = 5775
$ws1.Range("F7").Value = 9817
$ws1.Range("F20").Value = 634
$ws1.Range("F23").Value = 90
$ws1.Range("F26").Value = 2134
$ws1.Range("F28").Value = 365
$ws1.Range("F29").Value = 8064
$ws1.Range("F31").Value = 10
$ws1.Range("F41").Value = 1189
$ws1.Range("F42").Value = 1181
$ws1.Range("F45").Value = 2135
$ws1.Range("F48").Value = 1222

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 320
$ws4.Range("F6").Value = 5775
$ws4.Range("F18").Value = 634
$ws4.Range("F25").Value = 2134
$ws4.Range("F27").Value = 365
$ws4.Range("F28").Value = 8064
$ws4.Range("F38").Value = 1189
$ws4.Range("F39").Value = 1181
$ws4.Range("F44").Value = 2135
$ws4.Range("F49").Value = 1222

$wb.Save()
